$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 7 (Colombian Primera A / Tolima vs Junior FC Barranquilla)
# Its data now lives in row 6 with refreshed odds, so the whole row shifts up.
$ws.Range("A7:AO7").EntireRow.Delete()

# Refresh rows 2-6 with the latest league/match info and odds
# (column B "Date" is always 2025-12-16 and is left untouched)

# Row 2: Swiss Super League - St Gallen vs Sion
$ws.Cells.Item(2, 1).Value = "Swiss Super League"
$ws.Cells.Item(2, 3).Value = "16:30:00"
$ws.Cells.Item(2, 4).Value = "St Gallen"
$ws.Cells.Item(2, 5).Value = "Sion"
$ws.Cells.Item(2, 6).Value = 2.46
$ws.Cells.Item(2, 7).Value = 2.48
$ws.Cells.Item(2, 8).Value = 3.05
$ws.Cells.Item(2, 9).Value = 3.1
$ws.Cells.Item(2, 10).Value = 3.7
$ws.Cells.Item(2, 11).Value = 3.75
$ws.Cells.Item(2, 12).Value = 1.37
$ws.Cells.Item(2, 13).Value = 1.06
$ws.Cells.Item(2, 14).Value = 4.3
$ws.Cells.Item(2, 15).Value = 1.28
$ws.Cells.Item(2, 16).Value = 2.12
$ws.Cells.Item(2, 17).Value = 1.87
$ws.Cells.Item(2, 18).Value = 1.45
$ws.Cells.Item(2, 19).Value = 3.1
$ws.Cells.Item(2, 20).Value = 1.68
$ws.Cells.Item(2, 21).Value = 2.42
$ws.Cells.Item(2, 22).Value = 1.47
$ws.Cells.Item(2, 23).Value = 1.67
$ws.Cells.Item(2, 24).Value = 17
$ws.Cells.Item(2, 25).Value = 14
$ws.Cells.Item(2, 26).Value = 21
$ws.Cells.Item(2, 27).Value = 55
$ws.Cells.Item(2, 28).Value = 12
$ws.Cells.Item(2, 29).Value = 7.8
$ws.Cells.Item(2, 30).Value = 13.5
$ws.Cells.Item(2, 31).Value = 32
$ws.Cells.Item(2, 32).Value = 16.5
$ws.Cells.Item(2, 33).Value = 11.5
$ws.Cells.Item(2, 34).Value = 15.5
$ws.Cells.Item(2, 35).Value = 40
$ws.Cells.Item(2, 36).Value = 36
$ws.Cells.Item(2, 37).Value = 24
$ws.Cells.Item(2, 38).Value = 36
$ws.Cells.Item(2, 39).Value = 75
$ws.Cells.Item(2, 40).Value = 18
$ws.Cells.Item(2, 41).Value = 26

# Row 3: Swiss Super League - Winterthur vs Thun
$ws.Cells.Item(3, 1).Value = "Swiss Super League"
$ws.Cells.Item(3, 3).Value = "16:30:00"
$ws.Cells.Item(3, 4).Value = "Winterthur"
$ws.Cells.Item(3, 5).Value = "Thun"
$ws.Cells.Item(3, 6).Value = 4.6
$ws.Cells.Item(3, 7).Value = 4.8
$ws.Cells.Item(3, 8).Value = 1.72
$ws.Cells.Item(3, 9).Value = 1.74
$ws.Cells.Item(3, 10).Value = 4.6
$ws.Cells.Item(3, 11).Value = 4.8
$ws.Cells.Item(3, 12).Value = 1.28
$ws.Cells.Item(3, 13).Value = 1.03
$ws.Cells.Item(3, 14).Value = 6.4
$ws.Cells.Item(3, 15).Value = 1.17
$ws.Cells.Item(3, 16).Value = 2.84
$ws.Cells.Item(3, 17).Value = 1.52
$ws.Cells.Item(3, 18).Value = 1.74
$ws.Cells.Item(3, 19).Value = 2.3
$ws.Cells.Item(3, 20).Value = 1.59
$ws.Cells.Item(3, 21).Value = 2.6
$ws.Cells.Item(3, 22).Value = 2.34
$ws.Cells.Item(3, 23).Value = 1.26
$ws.Cells.Item(3, 24).Value = 28
$ws.Cells.Item(3, 25).Value = 13.5
$ws.Cells.Item(3, 26).Value = 13.5
$ws.Cells.Item(3, 27).Value = 18.5
$ws.Cells.Item(3, 28).Value = 26
$ws.Cells.Item(3, 29).Value = 11
$ws.Cells.Item(3, 30).Value = 10
$ws.Cells.Item(3, 31).Value = 15.5
$ws.Cells.Item(3, 32).Value = 42
$ws.Cells.Item(3, 33).Value = 18.5
$ws.Cells.Item(3, 34).Value = 16
$ws.Cells.Item(3, 35).Value = 24
$ws.Cells.Item(3, 36).Value = 95
$ws.Cells.Item(3, 37).Value = 48
$ws.Cells.Item(3, 38).Value = 44
$ws.Cells.Item(3, 39).Value = 65
$ws.Cells.Item(3, 40).Value = 36
$ws.Cells.Item(3, 41).Value = 6.8

# Row 4: English National League - Truro City vs Wealdstone
$ws.Cells.Item(4, 1).Value = "English National League"
$ws.Cells.Item(4, 3).Value = "16:45:00"
$ws.Cells.Item(4, 4).Value = "Truro City"
$ws.Cells.Item(4, 5).Value = "Wealdstone"
$ws.Cells.Item(4, 6).Value = 3.15
$ws.Cells.Item(4, 7).Value = 3.2
$ws.Cells.Item(4, 8).Value = 2.4
$ws.Cells.Item(4, 9).Value = 2.44
$ws.Cells.Item(4, 10).Value = 3.65
$ws.Cells.Item(4, 11).Value = 3.75
$ws.Cells.Item(4, 12).Value = 1.45
$ws.Cells.Item(4, 13).Value = 1.07
$ws.Cells.Item(4, 14).Value = 3.9
$ws.Cells.Item(4, 15).Value = 1.32
$ws.Cells.Item(4, 16).Value = 1.96
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(4, 18).Value = 1.38
$ws.Cells.Item(4, 19).Value = 3.55
$ws.Cells.Item(4, 20).Value = 1.74
$ws.Cells.Item(4, 21).Value = 2.2
$ws.Cells.Item(4, 22).Value = 1.7
$ws.Cells.Item(4, 23).Value = 1.45
$ws.Cells.Item(4, 24).Value = 15.5
$ws.Cells.Item(4, 25).Value = 10.5
$ws.Cells.Item(4, 26).Value = 15.5
$ws.Cells.Item(4, 27).Value = 34
$ws.Cells.Item(4, 28).Value = 13
$ws.Cells.Item(4, 29).Value = 8
$ws.Cells.Item(4, 30).Value = 11.5
$ws.Cells.Item(4, 31).Value = 25
$ws.Cells.Item(4, 32).Value = 21
$ws.Cells.Item(4, 33).Value = 13.5
$ws.Cells.Item(4, 34).Value = 17
$ws.Cells.Item(4, 35).Value = 38
$ws.Cells.Item(4, 36).Value = 55
$ws.Cells.Item(4, 37).Value = 36
$ws.Cells.Item(4, 38).Value = 48
$ws.Cells.Item(4, 39).Value = 140
$ws.Cells.Item(4, 40).Value = 36
$ws.Cells.Item(4, 41).Value = 21

# Row 5: Welsh Premiership - Cardiff Metropolitan vs Briton Ferry Llansawel
$ws.Cells.Item(5, 1).Value = "Welsh Premiership"
$ws.Cells.Item(5, 3).Value = "16:45:00"
$ws.Cells.Item(5, 4).Value = "Cardiff Metropolitan"
$ws.Cells.Item(5, 5).Value = "Briton Ferry Llansawel"
$ws.Cells.Item(5, 6).Value = 1.88
$ws.Cells.Item(5, 7).Value = 1.89
$ws.Cells.Item(5, 8).Value = 4.1
$ws.Cells.Item(5, 9).Value = 4.2
$ws.Cells.Item(5, 10).Value = 4.3
$ws.Cells.Item(5, 11).Value = 4.5
$ws.Cells.Item(5, 12).Value = 1.32
$ws.Cells.Item(5, 13).Value = 1.04
$ws.Cells.Item(5, 14).Value = 5.2
$ws.Cells.Item(5, 15).Value = 1.21
$ws.Cells.Item(5, 16).Value = 2.44
$ws.Cells.Item(5, 17).Value = 1.66
$ws.Cells.Item(5, 18).Value = 1.56
$ws.Cells.Item(5, 19).Value = 2.66
$ws.Cells.Item(5, 20).Value = 1.62
$ws.Cells.Item(5, 21).Value = 2.36
$ws.Cells.Item(5, 22).Value = 1.3
$ws.Cells.Item(5, 23).Value = 2.12
$ws.Cells.Item(5, 24).Value = 26
$ws.Cells.Item(5, 25).Value = 22
$ws.Cells.Item(5, 26).Value = 36
$ws.Cells.Item(5, 27).Value = 90
$ws.Cells.Item(5, 28).Value = 12.5
$ws.Cells.Item(5, 29).Value = 10.5
$ws.Cells.Item(5, 30).Value = 18
$ws.Cells.Item(5, 31).Value = 48
$ws.Cells.Item(5, 32).Value = 13.5
$ws.Cells.Item(5, 33).Value = 10.5
$ws.Cells.Item(5, 34).Value = 17.5
$ws.Cells.Item(5, 35).Value = 46
$ws.Cells.Item(5, 36).Value = 21
$ws.Cells.Item(5, 37).Value = 17
$ws.Cells.Item(5, 38).Value = 27
$ws.Cells.Item(5, 39).Value = 75
$ws.Cells.Item(5, 40).Value = 9
$ws.Cells.Item(5, 41).Value = 36

# Row 6: Colombian Primera A - Tolima vs Junior FC Barranquilla
$ws.Cells.Item(6, 1).Value = "Colombian Primera A"
$ws.Cells.Item(6, 3).Value = "21:30:00"
$ws.Cells.Item(6, 4).Value = "Tolima"
$ws.Cells.Item(6, 5).Value = "Junior FC Barranquilla"
$ws.Cells.Item(6, 6).Value = 1.76
$ws.Cells.Item(6, 7).Value = 1.8
$ws.Cells.Item(6, 8).Value = 5.3
$ws.Cells.Item(6, 9).Value = 5.8
$ws.Cells.Item(6, 10).Value = 3.7
$ws.Cells.Item(6, 11).Value = 3.95
$ws.Cells.Item(6, 12).Value = 1.46
$ws.Cells.Item(6, 13).Value = 1.08
$ws.Cells.Item(6, 14).Value = 3.25
$ws.Cells.Item(6, 15).Value = 1.4
$ws.Cells.Item(6, 16).Value = 1.78
$ws.Cells.Item(6, 17).Value = 2.22
$ws.Cells.Item(6, 18).Value = 1.28
$ws.Cells.Item(6, 19).Value = 4.2
$ws.Cells.Item(6, 20).Value = 2.04
$ws.Cells.Item(6, 21).Value = 1.83
$ws.Cells.Item(6, 22).Value = 1.2
$ws.Cells.Item(6, 23).Value = 2.24
$ws.Cells.Item(6, 24).Value = 13.5
$ws.Cells.Item(6, 25).Value = 16.5
$ws.Cells.Item(6, 26).Value = 42
$ws.Cells.Item(6, 27).Value = 1000
$ws.Cells.Item(6, 28).Value = 7.4
$ws.Cells.Item(6, 29).Value = 8.6
$ws.Cells.Item(6, 30).Value = 22
$ws.Cells.Item(6, 31).Value = 90
$ws.Cells.Item(6, 32).Value = 9.6
$ws.Cells.Item(6, 33).Value = 11
$ws.Cells.Item(6, 34).Value = 24
$ws.Cells.Item(6, 35).Value = 100
$ws.Cells.Item(6, 36).Value = 18.5
$ws.Cells.Item(6, 37).Value = 26
$ws.Cells.Item(6, 38).Value = 44
$ws.Cells.Item(6, 39).Value = 160
$ws.Cells.Item(6, 40).Value = 14.5
$ws.Cells.Item(6, 41).Value = 130
